$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 -> Balite
$ws.Range("B2").Value = "Balite"
$ws.Range("C2").Value = 14.8956
$ws.Range("D2").Value = 120.7855
$ws.Range("E2").Value = 5016
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 12

# Row 3 -> Balungao
$ws.Range("B3").Value = "Balungao"
$ws.Range("C3").Value = 14.9143
$ws.Range("D3").Value = 120.7622
$ws.Range("E3").Value = 5720
$ws.Range("F3").Value = 0
$ws.Range("G3").Value = 12

# Row 4 -> Bulusan
$ws.Range("B4").Value = "Bulusan"
$ws.Range("C4").Value = 14.9076
$ws.Range("D4").Value = 120.7455
$ws.Range("E4").Value = 2603
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 12

# Row 5 -> Frances
$ws.Range("B5").Value = "Frances"
$ws.Range("C5").Value = 14.9153
$ws.Range("D5").Value = 120.7532
$ws.Range("E5").Value = 6129
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = 12

# Row 6 (new) -> Gatbuca
$ws.Range("A6").Value = $true
$ws.Range("B6").Value = "Gatbuca"
$ws.Range("C6").Value = 14.9218
$ws.Range("D6").Value = 120.7685
$ws.Range("E6").Value = 6384
$ws.Range("F6").Value = 115
$ws.Range("G6").Value = 12

# Row 7 (new) -> Iba O'Este
$ws.Range("A7").Value = $true
$ws.Range("B7").Value = "Iba O'Este"
$ws.Range("C7").Value = 14.8919
$ws.Range("D7").Value = 120.7635
$ws.Range("E7").Value = 14085
$ws.Range("F7").Value = 601
$ws.Range("G7").Value = 12
